# Update "想去人数" (interested-people count) values in column F
# for the sheets "展览" and "全部类型", rows 5, 10, 20, 21.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F5").Value = 909
    $ws.Range("F10").Value = 4709
    $ws.Range("F20").Value = 3623
    $ws.Range("F21").Value = 266
}
